# Etat art II - Tests 50 articles
# Adds a new "SPRINT 4" sheet after "SPRINT 3", fills it with the sprint-4
# work log, updates the SPRINT 3 log with a final "Sprint 3 review" entry,
# and refreshes the SPRINTS roll-up formula to include SPRINT 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. SPRINT 3 sheet: add the last log row (row 16) and extend the total
# ---------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("SPRINT 3")

$sprint3.Range("B15").Copy() | Out-Null
$sprint3.Range("B16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$sprint3.Range("B16").Value = 43944
$sprint3.Range("C16").Value = "Sprint 3 review + rencontre TB"
$sprint3.Range("D16").Value = 1

$sprint3.Range("D17").Formula = "=SUM(D3:D16)"

# ---------------------------------------------------------------------
# 2. Create the new "SPRINT 4" sheet after "SPRINT 3"
# ---------------------------------------------------------------------
$sprint4 = $wb.Worksheets.Add($null, $sprint3)
$sprint4.Name = "SPRINT 4"

$sprint4.Range("A1").Formula = "=SUM(D3:D11)"

$sprint4.Range("B2").Value = "Date"
$sprint4.Range("C2").Value = "Quoi"
$sprint4.Range("D2").Value = "Temps (h)"

$rows = @(
    @(43951, "Test DB", 0.5),
    @(43951, "Manually create raw .txt files for articles", 2),
    @(43951, "Meeting with Zhan : front-end demo", 1),
    @(43953, "Raw text files - 50 articles", 5),
    @(43953, "PHP script loop test", 2),
    @(43956, "Raw text files - 50 articles", 2),
    @(43956, "Copyleaks - Test", 1),
    @(43956, "Unicheck - Tests 50", 1),
    @(43957, "Prepostseo - Tests 50, manually", 2),
    @(43957, "PlagiarismSearch - Tests 50 - script", 5)
)

$r = 3
foreach ($row in $rows) {
    $sprint4.Range("B$r").Value = $row[0]
    $sprint4.Range("C$r").Value = $row[1]
    $sprint4.Range("D$r").Value = $row[2]
    $r = $r + 1
}

# Date-format column B (rows 3-16) the same way as the other SPRINT sheets
$sprint3.Range("B3").Copy() | Out-Null
$sprint4.Range("B3:B16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$sprint4.Range("B3").Value = 43951

$sprint4.Range("D17").Formula = "=SUM(D3:D16)"

$sprint4.Columns.Item(3).ColumnWidth = 32.03

# ---------------------------------------------------------------------
# 3. SPRINTS sheet: include SPRINT 4 in the total-hours formula
# ---------------------------------------------------------------------
$sprints = $wb.Worksheets.Item("SPRINTS")
$sprints.Range("B1").Formula = "='SPRINT 0'!A1+'SPRINT 1'!A1+'SPRINT 2'!A1+'SPRINT 3'!A1+'SPRINT 4'!A1"
$sprints.Range("F9").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. SPRINT 2 sheet: selection moved
# ---------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("SPRINT 2")
$sprint2.Range("H13").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. SPRINT 3 sheet: selection moved, make SPRINT 4 the active tab
# ---------------------------------------------------------------------
$sprint3.Range("C19").Select() | Out-Null

$sprint4.Activate() | Out-Null
$sprint4.Range("B13").Select() | Out-Null
